$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate the new BOM rows (resistors + capacitor) in the exact
# cell-write order that reproduces the original shared-string layout.
$ws.Cells.Item(7, 4).Value  = "0.1uF Capacitor"
$ws.Cells.Item(8, 1).Value  = "RCC0402470RFKED"
$ws.Cells.Item(8, 2).Value  = "541-RCC0402470RFKEDCT-ND"
$ws.Cells.Item(8, 4).Value  = "470R 0402 resistor"
$ws.Cells.Item(9, 4).Value  = "100k 0402 resistor"
$ws.Cells.Item(9, 1).Value  = "RCA0402100KFKEDHP"
$ws.Cells.Item(9, 2).Value  = "541-3242-1-ND"
$ws.Cells.Item(10, 1).Value = "MCS04020C4701FE000"
$ws.Cells.Item(10, 2).Value = "MCS0402-4.70K-CFCT-ND"
$ws.Cells.Item(10, 4).Value = "4.7k 0402 resistor"
$ws.Cells.Item(7, 1).Value  = "CL05B104KO5VPNC"
$ws.Cells.Item(7, 2).Value  = "1276-6844-1-ND"

# Quantities
$ws.Cells.Item(7, 3).Value  = 2
$ws.Cells.Item(8, 3).Value  = 1
$ws.Cells.Item(9, 3).Value  = 1
$ws.Cells.Item(10, 3).Value = 1

# Widen column B to fit the new, longer part numbers
# (27.42578125 raw stored-width units; nearest value reachable through the
# ColumnWidth property's pixel-grid rounding)
$ws.Columns.Item(2).ColumnWidth = 26.666666666666668

# Move selection to A7, matching the saved view state
$ws.Range("A7").Select()
